$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 430
$ws.Range("I12").Value = 360
$ws.Range("K12").Value = 360
$ws.Range("M12").Value = -190

$ws.Range("H32").Value = 1387.1666
$ws.Range("J32").Value = 1387.1666
$ws.Range("L32").Value = 1387.1666
$ws.Range("N32").Value = -2039.1666

$ws.Range("H40").Value = 250001950
$ws.Range("J40").Value = 250001950
$ws.Range("L40").Value = 250001950
$ws.Range("N40").Value = -250002300

$ws.Range("H43").Value = 2396.8572
$ws.Range("J43").Value = 1815.6
$ws.Range("L43").Value = 1815.6
$ws.Range("N43").Value = -1953.6

$ws.Range("H51").Value = 7046.433
$ws.Range("J51").Value = 5614
$ws.Range("L51").Value = 5614
$ws.Range("N51").Value = -6582

$ws.Range("H141").Value = 26326114
$ws.Range("I141").Value = 35718700
$ws.Range("K141").Value = 107156100
$ws.Range("M141").Value = -107150920

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 6709.7
$ws.Range("I4").Value = 4498.5
$ws.Range("K4").Value = 4498.5
$ws.Range("M4").Value = -4382.5

$ws.Range("H60").Value = 68186.87
$ws.Range("I60").Value = 68186.87
$ws.Range("K60").Value = 68186.87
$ws.Range("M60").Value = -67453.87

$ws.Range("H61").Value = 54616428
$ws.Range("I61").Value = 87500870
$ws.Range("J61").Value = 2001318.6
$ws.Range("K61").Value = 87500870
$ws.Range("L61").Value = 2001318.6
$ws.Range("M61").Value = -87500658
$ws.Range("N61").Value = -2001742.6

$ws.Range("H74").Value = 4922.6
$ws.Range("I74").Value = 3446.6
$ws.Range("K74").Value = 3446.6
$ws.Range("M74").Value = -2572.6

$ws.Range("H77").Value = 4922.6
$ws.Range("I77").Value = 3446.6
$ws.Range("K77").Value = 17233
$ws.Range("M77").Value = -12865

$ws.Range("H102").Value = 2001.6
$ws.Range("J102").Value = 2322.6667
$ws.Range("L102").Value = 2322.6667
$ws.Range("N102").Value = -5566.6667

$ws.Range("H132").Value = 3708203.2
$ws.Range("I132").Value = 4679.56
$ws.Range("K132").Value = 14038.68
$ws.Range("M132").Value = -11508.68

$ws.Range("H136").Value = 54616428
$ws.Range("I136").Value = 87500870
$ws.Range("J136").Value = 2001318.6
$ws.Range("K136").Value = 262502610
$ws.Range("L136").Value = 6003955.800000001
$ws.Range("M136").Value = -262500060
$ws.Range("N136").Value = -6009055.800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 14288629
$ws.Range("I134").Value = 2804
$ws.Range("J134").Value = 25002998
$ws.Range("K134").Value = 8412
$ws.Range("L134").Value = 75008994
$ws.Range("M134").Value = -5877
$ws.Range("N134").Value = -75014064

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 58828190
$ws.Range("I31").Value = 71433016
$ws.Range("J31").Value = 5707.6665
$ws.Range("K31").Value = 71433016
$ws.Range("L31").Value = 5707.6665
$ws.Range("M31").Value = -71432721
$ws.Range("N31").Value = -6297.6665

$ws.Range("H34").Value = 58828190
$ws.Range("I34").Value = 71433016
$ws.Range("J34").Value = 5707.6665
$ws.Range("K34").Value = 71433016
$ws.Range("L34").Value = 5707.6665
$ws.Range("M34").Value = -71432814
$ws.Range("N34").Value = -6111.6665

$ws.Range("H58").Value = 3170
$ws.Range("I58").Value = 3016.5715
$ws.Range("K58").Value = 3016.5715
$ws.Range("M58").Value = -2813.5715

$ws.Range("H60").Value = 2944.5

$ws.Range("H134").Value = 4143.857
$ws.Range("I134").Value = 4143.857
$ws.Range("K134").Value = 12431.571
$ws.Range("M134").Value = -9896.571

$ws.Range("H136").Value = 3170
$ws.Range("I136").Value = 3016.5715
$ws.Range("K136").Value = 9049.7145
$ws.Range("M136").Value = -6499.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 98222.11
$ws.Range("J37").Value = 98222.11
$ws.Range("L37").Value = 294666.33
$ws.Range("N37").Value = -294890.33

$ws.Range("H113").Value = 1058.7693
$ws.Range("J113").Value = 1428.5
$ws.Range("L113").Value = 4285.5
$ws.Range("N113").Value = -8625.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 58.4
$ws.Range("J2").Value = 42.25
$ws.Range("L2").Value = 42.25
$ws.Range("N2").Value = -268.25

$ws.Range("H43").Value = 8993.333000000001
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H49").Value = 45000
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

$ws.Range("H80").Value = 60002770
$ws.Range("I80").Value = 2829.6667
$ws.Range("J80").Value = 85717030
$ws.Range("K80").Value = 2829.6667
$ws.Range("L80").Value = 85717030
$ws.Range("M80").Value = -1831.6667
$ws.Range("N80").Value = -85719026

$ws.Range("H83").Value = 60002770
$ws.Range("I83").Value = 2829.6667
$ws.Range("J83").Value = 85717030
$ws.Range("K83").Value = 14148.3335
$ws.Range("L83").Value = 428585150
$ws.Range("M83").Value = -9156.333500000001
$ws.Range("N83").Value = -428595134

$ws.Range("H122").Value = 3899.8333
$ws.Range("I122").Value = 3499.6667
$ws.Range("K122").Value = 10499.0001
$ws.Range("M122").Value = -8049.000100000001

$ws.Range("H132").Value = 12733329
$ws.Range("I132").Value = 4664.6665
$ws.Range("J132").Value = 15915496
$ws.Range("K132").Value = 13993.9995
$ws.Range("L132").Value = 47746488
$ws.Range("M132").Value = -11463.9995
$ws.Range("N132").Value = -47751548

$ws.Range("H135").Value = 108888.86
$ws.Range("J135").Value = 108888.86
$ws.Range("L135").Value = 108888.86
$ws.Range("N135").Value = -119028.86

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5610.8
$ws.Range("I40").Value = 5513.5
$ws.Range("K40").Value = 5513.5
$ws.Range("M40").Value = -5377.5

$ws.Range("H46").Value = 2277.1765
$ws.Range("I46").Value = 1858.4546
$ws.Range("J46").Value = 3044.8333
$ws.Range("K46").Value = 1858.4546
$ws.Range("L46").Value = 3044.8333
$ws.Range("M46").Value = -1670.4546
$ws.Range("N46").Value = -3420.8333

$ws.Range("H68").Value = 5720741.5
$ws.Range("I68").Value = 7480108.5
$ws.Range("K68").Value = 7480108.5
$ws.Range("M68").Value = -7479359.5

$ws.Range("H71").Value = 5720741.5
$ws.Range("I71").Value = 7480108.5
$ws.Range("K71").Value = 37400542.5
$ws.Range("M71").Value = -37396798.5

$ws.Range("H76").Value = 54999.5
$ws.Range("J76").Value = 54999.5
$ws.Range("L76").Value = 54999.5
$ws.Range("N76").Value = -55675.5

$ws.Range("H79").Value = 54999.5
$ws.Range("J79").Value = 54999.5
$ws.Range("L79").Value = 54999.5
$ws.Range("N79").Value = -57339.5

$ws.Range("H93").Value = 2319206.2
$ws.Range("I93").Value = 2245.2856
$ws.Range("J93").Value = 5562952
$ws.Range("K93").Value = 2245.2856
$ws.Range("L93").Value = 5562952
$ws.Range("M93").Value = -997.2856000000002
$ws.Range("N93").Value = -5565448

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 42000
$ws.Range("J54").Value = 42000
$ws.Range("L54").Value = 42000
$ws.Range("N54").Value = -43040

$ws.Range("H58").Value = 39745
$ws.Range("I58").Value = 49990
$ws.Range("J58").Value = 29500
$ws.Range("K58").Value = 49990
$ws.Range("L58").Value = 29500
$ws.Range("M58").Value = -49682
$ws.Range("N58").Value = -30116

$ws.Range("H100").Value = 618231.4
$ws.Range("I100").Value = 1613.1875
$ws.Range("K100").Value = 3226.375
$ws.Range("M100").Value = -2685.375

$ws.Range("H136").Value = 637962.8
$ws.Range("I136").Value = 13827
$ws.Range("K136").Value = 41481
$ws.Range("M136").Value = -38931
